$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New column E: width + header + data + "Bộ phận thực hiện"-style follow up
#    column ("Phản hồi từ bộ phận Sales") added to the right of the table.
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 23.6

# Give every new E cell the same thin border already used throughout the
# table (copy format from an existing plain-bordered cell so the engine
# reuses the existing border/xf definitions instead of inventing new ones).
$ws.Range("D13").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("E5").PasteSpecial(-4122)
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("E7").PasteSpecial(-4122)
$ws.Range("E8").PasteSpecial(-4122)
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("E11").PasteSpecial(-4122)
$ws.Range("E12").PasteSpecial(-4122)
$ws.Range("E13").PasteSpecial(-4122)
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("E15").PasteSpecial(-4122)

# Header cell E3
$ws.Range("E3").Value = "Phản hồi từ bộ phận Sales"
$ws.Range("E3").Font.Bold = $true
$ws.Range("E3").HorizontalAlignment = -4108
$ws.Range("E3").VerticalAlignment = -4108
$ws.Range("E3").WrapText = $true

# Data rows in column E
$ws.Range("E4").Value = "Admin PKD có trao đổi với chị Thảo kế toán để phản hồi về việc gửi link sau khi đã nộp thuế. "
$ws.Range("E4").WrapText = $true

$ws.Range("E6").Value = "Sau khi đã nộp thuế và phí hãng tàu, thông tin sẽ được hiển thị ở phần dashboard (đã update)"
$ws.Range("E6").WrapText = $true

$ws.Range("E7").Value = "Phần mềm đã update chức năng này"
$ws.Range("E7").WrapText = $true

$ws.Range("E8").Value = "Khi sale chốt hàng, admin phòng bán đã cập nhật lên hệ thống để trừ tồn. Trừ trường hợp PKD chưa bán được thì chưa có kế hoạch. PKD đảm bảo không phát sinh hạn lưu cont lưu bãi khi đến hạn lưu"
$ws.Range("E8").WrapText = $true

# ---------------------------------------------------------------------------
# 2. Make the header row (row 3) bold & centered both ways.
# ---------------------------------------------------------------------------
$ws.Range("A3:D3").Font.Bold = $true
$ws.Range("A3:D3").HorizontalAlignment = -4108
$ws.Range("A3:D3").VerticalAlignment = -4108
$ws.Rows.Item(3).RowHeight = 27.75

# ---------------------------------------------------------------------------
# 3. Row height tweaks (rows whose wrapped content height changed once the
#    new column / header styling was introduced).
# ---------------------------------------------------------------------------
$ws.Rows.Item(6).RowHeight = 60
$ws.Rows.Item(8).RowHeight = 135
$ws.Rows.Item(11).RowHeight = 205.5

# ---------------------------------------------------------------------------
# 4. Page setup (portrait, A4) and the selected/active cell & top-left view.
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1
$ws.PageSetup.PaperSize = 9

$ws.Range("D5").Select()
